$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2233.44
$ws.Range("I40").Value = 2133.6316
$ws.Range("J40").Value = 2549.5
$ws.Range("K40").Value = 2133.6316
$ws.Range("L40").Value = 2549.5
$ws.Range("M40").Value = -1958.6316
$ws.Range("N40").Value = -2899.5
$ws.Range("H53").Value = 520
$ws.Range("J53").Value = 535.8570999999999
$ws.Range("L53").Value = 535.8570999999999
$ws.Range("N53").Value = -1809.8571
$ws.Range("H80").Value = 159594.86
$ws.Range("I80").Value = 8226
$ws.Range("K80").Value = 24678
$ws.Range("M80").Value = -23680
$ws.Range("H83").Value = 159594.86
$ws.Range("I83").Value = 8226
$ws.Range("K83").Value = 74034
$ws.Range("M83").Value = -69042
$ws.Range("H98").Value = 2026.6296
$ws.Range("I98").Value = 786.35
$ws.Range("K98").Value = 786.35
$ws.Range("M98").Value = 711.65
$ws.Range("H101").Value = 7693735.5
$ws.Range("I101").Value = 20000676
$ws.Range("K101").Value = 60002028
$ws.Range("M101").Value = -60000406
$ws.Range("H118").Value = 1284.5
$ws.Range("I118").Value = 1284.5
$ws.Range("K118").Value = 3853.5
$ws.Range("M118").Value = -2196.5
$ws.Range("H122").Value = 2026.6296
$ws.Range("I122").Value = 786.35
$ws.Range("K122").Value = 2359.05
$ws.Range("M122").Value = 90.94999999999982
$ws.Range("H127").Value = 1711.2222
$ws.Range("I127").Value = 915.2857
$ws.Range("J127").Value = 4497
$ws.Range("K127").Value = 2745.8571
$ws.Range("L127").Value = 13491
$ws.Range("M127").Value = 2214.1429
$ws.Range("N127").Value = -23411
$ws.Range("H129").Value = 2426.2
$ws.Range("I129").Value = 816.25
$ws.Range("J129").Value = 3499.5
$ws.Range("K129").Value = 2448.75
$ws.Range("L129").Value = 10498.5
$ws.Range("M129").Value = 2551.25
$ws.Range("N129").Value = -20498.5
$ws.Range("H132").Value = 1788.1052
$ws.Range("I132").Value = 1645.5883
$ws.Range("K132").Value = 4936.7649
$ws.Range("M132").Value = -2406.7649
$ws.Range("H138").Value = 5597.4146
$ws.Range("J138").Value = 5569.3677
$ws.Range("L138").Value = 16708.1031
$ws.Range("N138").Value = -26988.1031
$ws.Range("H141").Value = 2598.8
$ws.Range("I141").Value = 2684.4285
$ws.Range("K141").Value = 8053.2855
$ws.Range("M141").Value = -2873.2855

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H32").Value = 17413.838
$ws.Range("I32").Value = 10911.167
$ws.Range("J32").Value = 39708.715
$ws.Range("K32").Value = 10911.167
$ws.Range("L32").Value = 39708.715
$ws.Range("M32").Value = -10624.167
$ws.Range("N32").Value = -40282.715
$ws.Range("H45").Value = 2009.8462
$ws.Range("I45").Value = 1837.6364
$ws.Range("K45").Value = 1837.6364
$ws.Range("M45").Value = -1460.6364
$ws.Range("H110").Value = 2616.4
$ws.Range("I110").Value = 2683.5557
$ws.Range("K110").Value = 2683.5557
$ws.Range("M110").Value = -638.5556999999999
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -19900
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 3055.7368
$ws.Range("I132").Value = 2947.7222
$ws.Range("K132").Value = 8843.1666
$ws.Range("M132").Value = -6313.1666

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 40000
$ws.Range("J6").Value = 40000
$ws.Range("L6").Value = 40000
$ws.Range("N6").Value = -40226
$ws.Range("H99").Value = 3193.5557
$ws.Range("J99").Value = 3344.4
$ws.Range("L99").Value = 3344.4
$ws.Range("N99").Value = -6340.4
$ws.Range("H105").Value = 2004.1
$ws.Range("I105").Value = 1434.8
$ws.Range("K105").Value = 1434.8
$ws.Range("M105").Value = 312.2

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 69998
$ws.Range("J68").Value = 69998
$ws.Range("L68").Value = 69998
$ws.Range("N68").Value = -71496
$ws.Range("H71").Value = 69998
$ws.Range("J71").Value = 69998
$ws.Range("L71").Value = 209994
$ws.Range("N71").Value = -217482
$ws.Range("H105").Value = 1961.4
$ws.Range("I105").Value = 1759.4286
$ws.Range("K105").Value = 1759.4286
$ws.Range("M105").Value = -12.42859999999996
$ws.Range("H107").Value = 552.6
$ws.Range("I107").Value = 200.09091
$ws.Range("K107").Value = 200.09091
$ws.Range("M107").Value = 1719.90909
$ws.Range("H134").Value = 3699
$ws.Range("I134").Value = 3699
$ws.Range("K134").Value = 11097
$ws.Range("M134").Value = -8562

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6181
$ws.Range("I5").Value = 9668.333000000001
$ws.Range("J5").Value = 950
$ws.Range("K5").Value = 29004.999
$ws.Range("L5").Value = 2850
$ws.Range("M5").Value = -28892.999
$ws.Range("N5").Value = -3074
$ws.Range("H38").Value = 2061.7693
$ws.Range("I38").Value = 1527.2727
$ws.Range("J38").Value = 5001.5
$ws.Range("K38").Value = 4581.8181
$ws.Range("L38").Value = 15004.5
$ws.Range("M38").Value = -4234.8181
$ws.Range("N38").Value = -15698.5
$ws.Range("H88").Value = 15089.875
$ws.Range("J88").Value = 15932.533
$ws.Range("L88").Value = 47797.599
$ws.Range("N88").Value = -48653.599
$ws.Range("H91").Value = 15089.875
$ws.Range("J91").Value = 15932.533
$ws.Range("L91").Value = 47797.599
$ws.Range("N91").Value = -50761.599
$ws.Range("H131").Value = 1963.3334
$ws.Range("J131").Value = 3000
$ws.Range("L131").Value = 9000
$ws.Range("N131").Value = -19080
$ws.Range("H135").Value = 6181
$ws.Range("I135").Value = 9668.333000000001
$ws.Range("J135").Value = 950
$ws.Range("K135").Value = 87014.997
$ws.Range("L135").Value = 8550
$ws.Range("M135").Value = -84479.997
$ws.Range("N135").Value = -13620

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1556.3334
$ws.Range("J102").Value = 1765
$ws.Range("L102").Value = 1765
$ws.Range("N102").Value = -5009
$ws.Range("H123").Value = 30242
$ws.Range("J123").Value = 30242
$ws.Range("L123").Value = 30242
$ws.Range("N123").Value = -35142
$ws.Range("H132").Value = 2071.625
$ws.Range("J132").Value = 2449.5
$ws.Range("L132").Value = 7348.5
$ws.Range("N132").Value = -12408.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5114.143
$ws.Range("I40").Value = 4599.6665
$ws.Range("K40").Value = 4599.6665
$ws.Range("M40").Value = -4463.6665
$ws.Range("H55").Value = 2534.889
$ws.Range("I55").Value = 1332.2858
$ws.Range("K55").Value = 1332.2858
$ws.Range("M55").Value = -1159.2858
$ws.Range("H122").Value = 8044.3335
$ws.Range("I122").Value = 8437.8125
$ws.Range("J122").Value = 7472
$ws.Range("K122").Value = 25313.4375
$ws.Range("L122").Value = 22416
$ws.Range("M122").Value = -22863.4375
$ws.Range("N122").Value = -27316
$ws.Range("H123").Value = 54999
$ws.Range("J123").Value = 54999
$ws.Range("L123").Value = 54999
$ws.Range("N123").Value = -64799

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I107").Value = 613
$ws.Range("J107").Value = 699.25
$ws.Range("K107").Value = 1839
$ws.Range("L107").Value = 2097.75
$ws.Range("M107").Value = 81
$ws.Range("N107").Value = -5937.75
